$d = $word.ActiveDocument

# 1) "[TITRE DU PROJET]" -> "TITRE_DU_PROJET"
$d.Content.Find.Execute("[TITRE DU PROJET]", $true, $false, $false, $false, $false, $true, 1, $false, "TITRE_DU_PROJET", 2) | Out-Null

# 2) "[AUTEUR(S)]" -> "AUTEUR(S)"
$d.Content.Find.Execute("[AUTEUR(S)]", $true, $false, $false, $false, $false, $true, 1, $false, "AUTEUR(S)", 2) | Out-Null

# 3) "Synopsis" (section title) -> "Synopsis intégral"
$d.Content.Find.Execute("Synopsis", $true, $false, $false, $false, $false, $true, 1, $false, "Synopsis intégral", 2) | Out-Null

# 4) Insert a new blank paragraph (style "Explicationdonne") right before the
#    "Duis tincidunt..." placeholder paragraph.
$findRange = $d.Content
$findRange.Find.Execute("Duis tincidunt, libero sit amet semper venenatis") | Out-Null
$insertionPoint = $d.Range($findRange.Start, $findRange.Start)
$insertionPoint.InsertParagraphBefore()

# Locate the "Duis tincidunt..." paragraph again and style the blank paragraph
# that now immediately precedes it.
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Duis tincidunt")) {
        $blankPara = $d.Paragraphs.Item($i - 1)
        $blankPara.Style = "Explicationdonne"
        break
    }
}

# 5) Prepend the new instruction sentence to the "Duis tincidunt..." paragraph.
$findRange2 = $d.Content
$findRange2.Find.Execute("Duis tincidunt, libero sit amet semper venenatis") | Out-Null
$prefixPoint = $d.Range($findRange2.Start, $findRange2.Start)
$prefixPoint.InsertBefore("LE SYNOPSIS DOIT PRÉSENTER L’HISTOIRE INTÉGRALEMENT, DE SON COMMENCEMENT JUSQU’À SA FIN (COMPRISE). ")
